$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (45177 -> 45178) for every data row from row 2 through row 125.
$ws.Range("C2:C125").Value = 45178
